$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, copying the
# existing header style (bold + border + centered) from O1.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# New data values (columns B..Q, i.e. cols 2..17) are identical across
# all data rows (2..25).
$values = @(2, 2, 2, 1, 1, 1, 2, 2, 2, 1, 2, 2, 2, 1, 2, 2)

for ($r = 2; $r -le 25; $r++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
